$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.750.48"
$ws.Range("E2").Value = "  +3.78%  "

$ws.Range("D3").Value = "1.772.89"
$ws.Range("E3").Value = "  +2.50%  "

$ws.Range("D4").Value = "'0.9966"
$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").Value = "'243.07"
$ws.Range("E5").Value = "  +0.65%  "

$ws.Range("D6").Value = "'0.9966"
$ws.Range("E6").Value = "  +0.25%  "

$ws.Range("D7").Value = "'0.4850"
$ws.Range("E7").Value = "  -1.48%  "

$ws.Range("D8").Value = "'0.2652"
$ws.Range("E8").Value = "  +1.74%  "

$ws.Range("D9").Value = "'0.06222"
$ws.Range("E9").Value = "  -0.07%  "

$ws.Range("D10").Value = "1.756.16"
$ws.Range("E10").Value = "  +1.69%  "

$ws.Range("D11").Value = "'16.26"
$ws.Range("E11").Value = "  +3.01%  "

$ws.Range("D12").Value = "'0.06960"
$ws.Range("E12").Value = "  -0.26%  "

$ws.Range("D13").Value = "'0.6139"
$ws.Range("E13").Value = "  +0.43%  "

$ws.Range("D14").Value = "'4.564"
$ws.Range("E14").Value = "  +1.55%  "

$ws.Range("D15").Value = "'78.45"
$ws.Range("E15").Value = "  +1.84%  "

$ws.Range("D16").Value = "'0.9974"
$ws.Range("E16").Value = "  +0.28%  "

$ws.Range("D17").Value = "27.726.93"
$ws.Range("E17").Value = "  +4.70%  "

$ws.Range("D18").Value = "'0.9961"
$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("D19").Value = "'0.000007122"
$ws.Range("E19").Value = "  -0.64%  "

$ws.Range("D20").Value = "'11.66"
$ws.Range("E20").Value = "  +2.53%  "

$ws.Range("D21").Value = "1.981.16"
$ws.Range("E21").Value = "  +1.87%  "

$ws.Range("D22").Value = "'4.555"
$ws.Range("E22").Value = "  +2.74%  "

$ws.Range("D23").Value = "'8.532"
$ws.Range("E23").Value = "  +0.42%  "

$ws.Range("D24").Value = "'5.176"
$ws.Range("E24").Value = "  +0.88%  "

$ws.Range("D25").Value = "'141.49"
$ws.Range("E25").Value = "  +2.64%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'15.51"
$ws.Range("E26").Value = "  +1.24%  "

$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'1.873"
$ws.Range("E27").Value = "  +6.60%  "

$ws.Range("D28").Value = "'110.29"
$ws.Range("E28").Value = "  +3.69%  "

$ws.Range("D29").Value = "'1.383"
$ws.Range("E29").Value = "  -2.27%  "

$ws.Range("D30").Value = "'4.006"
$ws.Range("E30").Value = "  +1.91%  "

$ws.Range("D31").Value = "'0.08200"
$ws.Range("E31").Value = "  +2.95%  "

$ws.Range("D32").Value = "'3.728"
$ws.Range("E32").Value = "  +1.97%  "

$ws.Range("D33").Value = "'0.04677"
$ws.Range("E33").Value = "  +3.61%  "

$ws.Range("D34").Value = "'1.042"
$ws.Range("E34").Value = "  +3.71%  "

$ws.Range("D35").Value = "'2.605"
$ws.Range("E35").Value = "  -0.56%  "

$ws.Range("D36").Value = "'0.6289"
$ws.Range("E36").Value = "  +0.71%  "

$ws.Range("D37").Value = "'0.9327"
$ws.Range("E37").Value = "  -1.98%  "

$ws.Range("D38").Value = "'2.590"
$ws.Range("E38").Value = "  +7.22%  "

$ws.Range("D39").Value = "'2.038"
$ws.Range("E39").Value = "  +0.44%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01518"
$ws.Range("E40").Value = "  +1.25%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'5.755"
$ws.Range("E41").Value = "  +4.61%  "

$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'0.9969"
$ws.Range("E42").Value = "  +0.26%  "

$ws.Range("D43").Value = "'99.77"
$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("D44").Value = "'0.3898"
$ws.Range("E44").Value = "  +0.82%  "

$ws.Range("D45").Value = "'6.949"
$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").Value = "'0.1174"
$ws.Range("E46").Value = "  +0.84%  "

$ws.Range("D47").Value = "'0.05389"
$ws.Range("E47").Value = "  +0.24%  "

$ws.Range("D48").Value = "'7.901"
$ws.Range("E48").Value = "  +1.17%  "

$ws.Range("D49").Value = "'30.21"
$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("D50").Value = "'1.265"
$ws.Range("E50").Value = "  +3.40%  "

$ws.Range("D51").Value = "'52.05"
$ws.Range("E51").Value = "  +0.79%  "
